$wb = $excel.ActiveWorkbook

# Overview sheet: update status for rows 3 and 4 (files 2f1123f1-... and 6b9cff2a-...)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "In Translation"
$wsOverview.Range("C3").Value = "In Translation"
$wsOverview.Range("B4").Value = "In Translation"
$wsOverview.Range("C4").Value = "In Translation"

# zh-cn sheet: update status for rows 3 and 4
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = "In Translation"
$wsZhCn.Range("B4").Value = "In Translation"

# de-de sheet: update status for rows 3 and 4
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = "In Translation"
$wsDeDe.Range("B4").Value = "In Translation"
